$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 2724.6348400000002
$ws.Range("C7").Value = 5406.0033199999998
$ws.Range("D7").Value = 5139.5067499999996
$ws.Range("E7").Value = 8878.2346699999998
$ws.Range("F7").Value = 4469.6042699999998
$ws.Range("G7").Value = 11283.687599999999

$ws.Range("G8").Select()
